$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update store_attr_1_value / store_attr_2_value / store_attr_3_value
# text for each CC_* profile row (Tokyo POC KPI change: add shared
# Priority_POC / Other_POC markers + numbered scene codes).
$ws.Range("B2").Value = "Priority_POC_SM_L, Other_POC_SM_L, Priority_POC,Other_POC"
$ws.Range("C2").Value = "Priority_POC_SM_L, Priority_POC,Other_POC"
$ws.Range("D2").Value = "Event Space, Checkout, Beverage end, Deli section, Liquor, 0003-Event space, 0004-Checkout,0005-End, 0006-Deli,0007-Liquor"
$ws.Range("B3").Value = "Priority_POC_SM_M, Other_POC_SM_M, Priority_POC,Other_POC"
$ws.Range("C3").Value = "Priority_POC_SM_M, Priority_POC,Other_POC"
$ws.Range("D3").Value = "Event Space, Checkout, Deli section, Liquor, 0003-Event space,0004-Checkout,0006-Deli,0007-Liquor"
$ws.Range("B4").Value = "Priority_POC_SM_S, Other_POC_SM_S, Priority_POC,Other_POC"
$ws.Range("C4").Value = "Priority_POC_SM_S, Priority_POC,Other_POC"
$ws.Range("D4").Value = "Event Space, Checkout, Deli section, 0003-Event space,0004-Checkout,0006-Deli"
$ws.Range("B5").Value = "Priority_POC_Drug_L, Other_POC_Drug_L, Priority_POC,Other_POC"
$ws.Range("C5").Value = "Priority_POC_Drug_L, Priority_POC,Other_POC"
$ws.Range("D5").Value = "Event Space, Checkout, Beverage end, Deli section, 0003-Event space,0004-Checkout,0005-End,0006-Deli"
$ws.Range("B6").Value = "Priority_POC_Drug_M, Other_POC_Drug_M, Priority_POC,Other_POC"
$ws.Range("C6").Value = "Priority_POC_Drug_M, Priority_POC,Other_POC"
$ws.Range("D6").Value = "Event Space, Checkout, Store front, 0003-Event space,0004-Checkout, 0009-Store Front"
$ws.Range("B7").Value = "Priority_POC_Drug_S, Other_POC_Drug_S, Priority_POC,Other_POC"
$ws.Range("C7").Value = "Priority_POC_Drug_S, Priority_POC,Other_POC"
$ws.Range("D7").Value = "Checkout, Store front,0004-Checkout, 0009-Store Front"
$ws.Range("B8").Value = "Priority_POC_Discounter_L, Other_POC_Discounter_L, Priority_POC,Other_POC"
$ws.Range("C8").Value = "Priority_POC_Discounter_L, Priority_POC,Other_POC"
$ws.Range("D8").Value = "Event Space, Checkout, Beverage end, Liquor,0003-Event space,0004-Checkout,0005-End,0007-Liquor"
$ws.Range("B9").Value = "Priority_POC_Discounter_M, Other_POC_Discounter_M, Priority_POC,Other_POC"
$ws.Range("C9").Value = "Priority_POC_Discounter_M, Priority_POC,Other_POC"
$ws.Range("D9").Value = "Event Space, Checkout, Beverage end,0003-Event space,0004-Checkout,0005-End"
$ws.Range("B10").Value = "Priority_POC_Discounter_S, Other_POC_Discounter_S, Priority_POC,Other_POC"
$ws.Range("C10").Value = "Priority_POC_Discounter_S, Priority_POC,Other_POC"
$ws.Range("D10").Value = "Event Space, Checkout,0003-Event space,0004-Checkout"

# Widen columns B:D (and nudge the small KPI columns) to fit the new text.
$ws.Columns.Item(2).ColumnWidth = 45.01282051282046
$ws.Columns.Item(3).ColumnWidth = 48.97638326585697
$ws.Columns.Item(4).ColumnWidth = 92.46626180836707
$ws.Columns.Item(7).ColumnWidth = 12.235492577597865
$ws.Columns.Item(9).ColumnWidth = 12.235492577597865
$ws.Columns.Item(10).ColumnWidth = 12.340755735492566
$ws.Columns.Item(11).ColumnWidth = 12.235492577597865
$ws.Columns.Item(12).ColumnWidth = 12.340755735492566
$ws.Columns.Item(13).ColumnWidth = 12.235492577597865
$ws.Columns.Item(14).ColumnWidth = 12.340755735492566
$ws.Columns.Item(15).ColumnWidth = 12.235492577597865
$ws.Columns.Item(16).ColumnWidth = 12.340755735492566
$ws.Columns.Item(17).ColumnWidth = 12.235492577597865

# Move the saved selection/active cell (was M32).
$ws.Range("D25").Select() | Out-Null
